# The pre-seeded $wb can come back null in this host; $excel.ActiveWorkbook
# is the reliable handle, so re-resolve everything from there.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: column F "Predicted Eastern Elongation (UT)" ---
# F5:F18 are currently blank (General format). Give them the same
# number format/style the rest of column F already uses (carried by F3)
# BEFORE putting formulas in, so Excel doesn't auto-apply its own generic
# date format to the newly-populated cells.
$ws2.Range("F3").Copy()
$ws2.Range("F5:F18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# F4 already existed; its DATE(...) formula used day-of-month argument 1.
# Change it to 0 (i.e. last day of the *previous* month == one day earlier).
$ws2.Range("F4").Formula = '=DATE(YEAR(TODAY()), MONTH(1 & $A$2), 0) + INT(D4) + (MOD(D4, 1) * 24/24)'

# Fill the same formula (relative refs auto-adjust) down through F18, which
# previously had no formula in column F at all.
$ws2.Range("F5:F18").Formula = '=DATE(YEAR(TODAY()), MONTH(1 & $A$2), 0) + INT(D5) + (MOD(D5, 1) * 24/24)'

# --- Selections (cursor position only, no data impact) ---
$ws1.Activate()
$ws1.Range("F4").Select()

$ws2.Activate()
$ws2.Range("G1").Select()

# Restore Sheet1 as the active/visible tab.
$ws1.Activate()
